$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 43.5
$ws.Cells.Item(6, 9).Value = 43.5
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 130.5
$ws.Cells.Item(6, 12).Value = 0
$ws.Cells.Item(6, 13).Value = -18.5
$ws.Cells.Item(6, 14).ClearContents()

$ws.Cells.Item(17, 8).Value = 879.087
$ws.Cells.Item(17, 10).Value = 879.087
$ws.Cells.Item(17, 12).Value = 2637.261
$ws.Cells.Item(17, 14).Value = -2973.261

$ws.Cells.Item(33, 8).Value = 126.53333
$ws.Cells.Item(33, 9).Value = 114.53846
$ws.Cells.Item(33, 11).Value = 114.53846
$ws.Cells.Item(33, 13).Value = 114.46154

$ws.Cells.Item(86, 8).Value = 3746.611
$ws.Cells.Item(86, 9).Value = 3184.4443
$ws.Cells.Item(86, 11).Value = 3184.4443
$ws.Cells.Item(86, 13).Value = -2061.4443

$ws.Cells.Item(88, 8).Value = 1457599.8
$ws.Cells.Item(88, 10).Value = 264498.75
$ws.Cells.Item(88, 12).Value = 264498.75
$ws.Cells.Item(88, 14).Value = -265310.75

$ws.Cells.Item(89, 8).Value = 3746.611
$ws.Cells.Item(89, 9).Value = 3184.4443
$ws.Cells.Item(89, 11).Value = 15922.2215
$ws.Cells.Item(89, 13).Value = -10306.2215

$ws.Cells.Item(91, 8).Value = 1457599.8
$ws.Cells.Item(91, 10).Value = 264498.75
$ws.Cells.Item(91, 12).Value = 264498.75
$ws.Cells.Item(91, 14).Value = -267306.75

$ws.Cells.Item(98, 8).Value = 618.8077
$ws.Cells.Item(98, 9).Value = 623.6
$ws.Cells.Item(98, 11).Value = 623.6
$ws.Cells.Item(98, 13).Value = 874.4

$ws.Cells.Item(122, 8).Value = 618.8077
$ws.Cells.Item(122, 9).Value = 623.6
$ws.Cells.Item(122, 11).Value = 1870.8
$ws.Cells.Item(122, 13).Value = 579.1999999999998

$ws.Cells.Item(132, 8).Value = 1952.4584
$ws.Cells.Item(132, 9).Value = 1874.3478
$ws.Cells.Item(132, 10).Value = 3749
$ws.Cells.Item(132, 11).Value = 5623.0434
$ws.Cells.Item(132, 12).Value = 11247
$ws.Cells.Item(132, 13).Value = -3093.0434
$ws.Cells.Item(132, 14).Value = -16307

$ws.Cells.Item(137, 8).Value = 485777.56
$ws.Cells.Item(137, 9).Value = 1591.65
$ws.Cells.Item(137, 10).Value = 1454149.4
$ws.Cells.Item(137, 11).Value = 4774.950000000001
$ws.Cells.Item(137, 12).Value = 4362448.199999999
$ws.Cells.Item(137, 13).Value = -2224.950000000001
$ws.Cells.Item(137, 14).Value = -4367548.199999999

$ws.Cells.Item(138, 8).Value = 1541.5264
$ws.Cells.Item(138, 10).Value = 2039.56
$ws.Cells.Item(138, 12).Value = 6118.68
$ws.Cells.Item(138, 14).Value = -16398.68

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(39, 8).Value = 15819.833
$ws.Cells.Item(39, 9).Value = 8340
$ws.Cells.Item(39, 11).Value = 8340
$ws.Cells.Item(39, 13).Value = -7820

$ws.Cells.Item(61, 8).Value = 32287.182
$ws.Cells.Item(61, 9).Value = 1804.8214
$ws.Cells.Item(61, 11).Value = 1804.8214
$ws.Cells.Item(61, 13).Value = -1592.8214

$ws.Cells.Item(74, 8).Value = 26755.525
$ws.Cells.Item(74, 9).Value = 31391.637
$ws.Cells.Item(74, 11).Value = 31391.637
$ws.Cells.Item(74, 13).Value = -30517.637

$ws.Cells.Item(77, 8).Value = 26755.525
$ws.Cells.Item(77, 9).Value = 31391.637
$ws.Cells.Item(77, 11).Value = 156958.185
$ws.Cells.Item(77, 13).Value = -152590.185

$ws.Cells.Item(110, 8).Value = 2165.8572
$ws.Cells.Item(110, 9).Value = 1632.6666
$ws.Cells.Item(110, 11).Value = 1632.6666
$ws.Cells.Item(110, 13).Value = 412.3334

$ws.Cells.Item(122, 8).Value = 4002
$ws.Cells.Item(122, 10).Value = 4998
$ws.Cells.Item(122, 12).Value = 14994
$ws.Cells.Item(122, 14).Value = -19894

$ws.Cells.Item(132, 8).Value = 2190.9355
$ws.Cells.Item(132, 9).Value = 1780.6538
$ws.Cells.Item(132, 11).Value = 5341.9614
$ws.Cells.Item(132, 13).Value = -2811.9614

$ws.Cells.Item(136, 8).Value = 32287.182
$ws.Cells.Item(136, 9).Value = 1804.8214
$ws.Cells.Item(136, 11).Value = 5414.4642
$ws.Cells.Item(136, 13).Value = -2864.4642

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 402247.38
$ws.Cells.Item(20, 9).Value = 642837.8
$ws.Cells.Item(20, 11).Value = 642837.8
$ws.Cells.Item(20, 13).Value = -642590.8

$ws.Cells.Item(86, 8).Value = 3339.4167
$ws.Cells.Item(86, 9).Value = 3580.3333
$ws.Cells.Item(86, 10).Value = 2937.889
$ws.Cells.Item(86, 11).Value = 3580.3333
$ws.Cells.Item(86, 12).Value = 2937.889
$ws.Cells.Item(86, 13).Value = -2457.3333
$ws.Cells.Item(86, 14).Value = -5183.889

$ws.Cells.Item(89, 8).Value = 3339.4167
$ws.Cells.Item(89, 9).Value = 3580.3333
$ws.Cells.Item(89, 10).Value = 2937.889
$ws.Cells.Item(89, 11).Value = 17901.6665
$ws.Cells.Item(89, 12).Value = 14689.445
$ws.Cells.Item(89, 13).Value = -12285.6665
$ws.Cells.Item(89, 14).Value = -25921.445

$ws.Cells.Item(107, 8).Value = 6669598
$ws.Cells.Item(107, 9).Value = 9093833
$ws.Cells.Item(107, 11).Value = 9093833
$ws.Cells.Item(107, 13).Value = -9091913

$ws.Cells.Item(134, 8).Value = 4655.1665
$ws.Cells.Item(134, 9).Value = 3201.875
$ws.Cells.Item(134, 11).Value = 9605.625
$ws.Cells.Item(134, 13).Value = -7070.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1664.0769
$ws.Cells.Item(16, 9).Value = 1077.25
$ws.Cells.Item(16, 11).Value = 1077.25
$ws.Cells.Item(16, 13).Value = -790.25

$ws.Cells.Item(31, 8).Value = 4841.6665
$ws.Cells.Item(31, 9).Value = 3991.25
$ws.Cells.Item(31, 10).Value = 5522
$ws.Cells.Item(31, 11).Value = 3991.25
$ws.Cells.Item(31, 12).Value = 5522
$ws.Cells.Item(31, 13).Value = -3696.25
$ws.Cells.Item(31, 14).Value = -6112

$ws.Cells.Item(34, 8).Value = 4841.6665
$ws.Cells.Item(34, 9).Value = 3991.25
$ws.Cells.Item(34, 10).Value = 5522
$ws.Cells.Item(34, 11).Value = 3991.25
$ws.Cells.Item(34, 12).Value = 5522
$ws.Cells.Item(34, 13).Value = -3789.25
$ws.Cells.Item(34, 14).Value = -5926

$ws.Cells.Item(35, 8).Value = 3180.5557
$ws.Cells.Item(35, 9).Value = 3870.8333
$ws.Cells.Item(35, 10).Value = 1800
$ws.Cells.Item(35, 11).Value = 3870.8333
$ws.Cells.Item(35, 12).Value = 1800
$ws.Cells.Item(35, 13).Value = -3576.8333
$ws.Cells.Item(35, 14).Value = -2388

$ws.Cells.Item(58, 8).Value = 2009.3334
$ws.Cells.Item(58, 9).Value = 2000
$ws.Cells.Item(58, 11).Value = 2000
$ws.Cells.Item(58, 13).Value = -1797

$ws.Cells.Item(86, 8).Value = 8932846
$ws.Cells.Item(86, 10).Value = 7250
$ws.Cells.Item(86, 12).Value = 7250
$ws.Cells.Item(86, 14).Value = -9496

$ws.Cells.Item(89, 8).Value = 8932846
$ws.Cells.Item(89, 10).Value = 7250
$ws.Cells.Item(89, 12).Value = 36250
$ws.Cells.Item(89, 14).Value = -47482

$ws.Cells.Item(107, 8).Value = 1241.4828
$ws.Cells.Item(107, 9).Value = 1230.625
$ws.Cells.Item(107, 11).Value = 1230.625
$ws.Cells.Item(107, 13).Value = 689.375

$ws.Cells.Item(113, 8).Value = 1664.0769
$ws.Cells.Item(113, 9).Value = 1077.25
$ws.Cells.Item(113, 11).Value = 1077.25
$ws.Cells.Item(113, 13).Value = 1092.75

$ws.Cells.Item(132, 8).Value = 3440469.5
$ws.Cells.Item(132, 9).Value = 2842998
$ws.Cells.Item(132, 10).Value = 13000014
$ws.Cells.Item(132, 11).Value = 8528994
$ws.Cells.Item(132, 12).Value = 39000042
$ws.Cells.Item(132, 13).Value = -8526464
$ws.Cells.Item(132, 14).Value = -39005102

$ws.Cells.Item(136, 8).Value = 2009.3334
$ws.Cells.Item(136, 9).Value = 2000
$ws.Cells.Item(136, 11).Value = 6000
$ws.Cells.Item(136, 13).Value = -3450

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 85428.914
$ws.Cells.Item(23, 10).Value = 93016.27
$ws.Cells.Item(23, 12).Value = 279048.81
$ws.Cells.Item(23, 14).Value = -279518.81

$ws.Cells.Item(130, 8).Value = 11264.25

$ws.Cells.Item(136, 8).Value = 2748.3572
$ws.Cells.Item(136, 9).Value = 1553.5
$ws.Cells.Item(136, 10).Value = 4341.5
$ws.Cells.Item(136, 11).Value = 4660.5
$ws.Cells.Item(136, 12).Value = 13024.5
$ws.Cells.Item(136, 13).Value = 439.5
$ws.Cells.Item(136, 14).Value = -23224.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 939617.4
$ws.Cells.Item(11, 9).Value = 399434.44
$ws.Cells.Item(11, 10).Value = 1884937.5
$ws.Cells.Item(11, 11).Value = 399434.44
$ws.Cells.Item(11, 12).Value = 1884937.5
$ws.Cells.Item(11, 13).Value = -399295.44
$ws.Cells.Item(11, 14).Value = -1885215.5

$ws.Cells.Item(70, 8).Value = 39216.562
$ws.Cells.Item(70, 9).Value = 72052
$ws.Cells.Item(70, 10).Value = 6381.125
$ws.Cells.Item(70, 11).Value = 72052
$ws.Cells.Item(70, 12).Value = 6381.125
$ws.Cells.Item(70, 13).Value = -71782
$ws.Cells.Item(70, 14).Value = -6921.125

$ws.Cells.Item(73, 8).Value = 39216.562
$ws.Cells.Item(73, 9).Value = 72052
$ws.Cells.Item(73, 10).Value = 6381.125
$ws.Cells.Item(73, 11).Value = 72052
$ws.Cells.Item(73, 12).Value = 6381.125
$ws.Cells.Item(73, 13).Value = -71116
$ws.Cells.Item(73, 14).Value = -8253.125

$ws.Cells.Item(80, 8).Value = 3516.5
$ws.Cells.Item(80, 10).Value = 3175
$ws.Cells.Item(80, 12).Value = 3175
$ws.Cells.Item(80, 14).Value = -5171

$ws.Cells.Item(83, 8).Value = 3516.5
$ws.Cells.Item(83, 10).Value = 3175
$ws.Cells.Item(83, 12).Value = 15875
$ws.Cells.Item(83, 14).Value = -25859

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 239.07692
$ws.Cells.Item(61, 9).Value = 250.09091
$ws.Cells.Item(61, 11).Value = 250.09091
$ws.Cells.Item(61, 13).Value = -48.09091000000001

$ws.Cells.Item(113, 8).Value = 239.07692
$ws.Cells.Item(113, 9).Value = 250.09091
$ws.Cells.Item(113, 11).Value = 250.09091
$ws.Cells.Item(113, 13).Value = 1919.90909

$ws.Cells.Item(136, 8).Value = 2020.0476
$ws.Cells.Item(136, 9).Value = 1609.7693
$ws.Cells.Item(136, 11).Value = 4829.3079
$ws.Cells.Item(136, 13).Value = -2279.3079

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2274.875
$ws.Cells.Item(81, 9).Value = 899.8182
$ws.Cells.Item(81, 10).Value = 5300
$ws.Cells.Item(81, 11).Value = 1799.6364
$ws.Cells.Item(81, 12).Value = 10600
$ws.Cells.Item(81, 13).Value = -738.6364000000001
$ws.Cells.Item(81, 14).Value = -12722

$ws.Cells.Item(84, 8).Value = 2274.875
$ws.Cells.Item(84, 9).Value = 899.8182
$ws.Cells.Item(84, 10).Value = 5300
$ws.Cells.Item(84, 11).Value = 8998.182000000001
$ws.Cells.Item(84, 12).Value = 53000
$ws.Cells.Item(84, 13).Value = -3694.182000000001
$ws.Cells.Item(84, 14).Value = -63608

$ws.Cells.Item(114, 8).Value = 75000
$ws.Cells.Item(114, 10).Value = 75000
$ws.Cells.Item(114, 12).Value = 75000
$ws.Cells.Item(114, 14).Value = -83678

